$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.232.92'
$ws.Range("E2").Value = '  +0.91%  '

$ws.Range("D3").Value = '1.853.37'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("E5").Value = '  +0.92%  '

$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4650'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3715'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("E9").Value = '  -0.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8908'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07867'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("D13").Value = '1.831.42'
$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.405'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.33%  '

$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008925'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.95%  '

$ws.Range("E19").Value = '  -0.26%  '

$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("D21").Value = '27.262.34'
$ws.Range("E21").Value = '  +0.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.088'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").Value = '2.083.89'
$ws.Range("E24").Value = '  +0.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.963'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.57%  '

$ws.Range("E27").Value = '  -0.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.041'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.040'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08833'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.139'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7687'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.167'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.524'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.715'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.110'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.18%  '

$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.947'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.050'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.97%  '

$ws.Range("E42").Value = '  -0.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1626'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.490'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4796'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.71%  '

$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.644'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06202'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.91%  '
